$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Km initiali (B12)
$ws.Range("B12").Value = 52375

# Daily travel rows (A=day already set; update B/C/D)
$ws.Range("B15").Value = 356
$ws.Range("C15").Value = "Cluj-Baia-Mare"
$ws.Range("D15").Value = "Interes Serviciu"
$ws.Range("B17").Value = 152
$ws.Range("C17").Value = "Cluj-Cmp. Turzii"
$ws.Range("D17").Value = "Interes Serviciu"
$ws.Range("B19").Value = 356
$ws.Range("C19").Value = "Cluj-Baia-Mare"
$ws.Range("D19").Value = "Interes Serviciu"
$ws.Range("B21").Value = 152
$ws.Range("C21").Value = "Cluj-Cmp. Turzii"
$ws.Range("D21").Value = "Interes Serviciu"
$ws.Range("B23").Value = 85
$ws.Range("C23").Value = "Cluj-Apahida"
$ws.Range("D23").Value = "Interes Serviciu"
$ws.Range("B25").Value = 30
$ws.Range("C25").Value = "Acasa-Birou"
$ws.Range("D25").Value = " "
$ws.Range("B27").Value = 257
$ws.Range("C27").Value = "Cluj-Bistrita"
$ws.Range("D27").Value = "Interes Serviciu"
$ws.Range("B29").Value = 30
$ws.Range("C29").Value = "Acasa-Birou"
$ws.Range("D29").Value = " "
$ws.Range("B31").Value = 30
$ws.Range("C31").Value = "Acasa-Birou"
$ws.Range("D31").Value = " "
$ws.Range("B33").Value = 152
$ws.Range("C33").Value = "Cluj-Cmp. Turzii"
$ws.Range("D33").Value = "Interes Serviciu"
$ws.Range("B35").Value = 101
$ws.Range("C35").Value = "Cluj-Dej"
$ws.Range("D35").Value = "Interes Serviciu"
$ws.Range("B37").Value = 85
$ws.Range("C37").Value = "Cluj-Apahida"
$ws.Range("D37").Value = "Interes Serviciu"
$ws.Range("B39").Value = 421
$ws.Range("C39").Value = "Cluj-Satu-Mare"
$ws.Range("D39").Value = "Interes Serviciu"
$ws.Range("B41").Value = 156
$ws.Range("C41").Value = "Cluj-Zalau"
$ws.Range("D41").Value = "Interes Serviciu"
$ws.Range("B43").Value = 257
$ws.Range("C43").Value = "Cluj-Bistrita"
$ws.Range("D43").Value = "Interes Serviciu"
$ws.Range("B45").Value = 47
$ws.Range("C45").Value = "Cluj-Cluj"
$ws.Range("D45").Value = "Interes Serviciu"
$ws.Range("B47").Value = 356
$ws.Range("C47").Value = "Cluj-Baia-Mare"
$ws.Range("D47").Value = "Interes Serviciu"
$ws.Range("B49").Value = 30
$ws.Range("C49").Value = "Acasa-Birou"
$ws.Range("D49").Value = " "
$ws.Range("B51").Value = 30
$ws.Range("C51").Value = "Acasa-Birou"
$ws.Range("D51").Value = " "
$ws.Range("B53").Value = 121
$ws.Range("C53").Value = "Cluj-Turda"
$ws.Range("D53").Value = "Interes Serviciu"
$ws.Range("B55").Value = 356
$ws.Range("C55").Value = "Cluj-Baia-Mare"
$ws.Range("D55").Value = "Interes Serviciu"
$ws.Range("B57").Value = 101
$ws.Range("C57").Value = "Cluj-Dej"
$ws.Range("D57").Value = "Interes Serviciu"
$ws.Range("B59").Value = 421
$ws.Range("C59").Value = "Cluj-Satu-Mare"
$ws.Range("D59").Value = "Interes Serviciu"
$ws.Range("B61").Value = 101
$ws.Range("C61").Value = "Cluj-Dej"
$ws.Range("D61").Value = "Interes Serviciu"
$ws.Range("B63").Value = 156
$ws.Range("C63").Value = "Cluj-Zalau"
$ws.Range("D63").Value = "Interes Serviciu"
$ws.Range("B65").Value = 30
$ws.Range("C65").Value = "Acasa-Birou"
$ws.Range("D65").Value = " "
$ws.Range("B67").Value = 356
$ws.Range("C67").Value = "Cluj-Baia-Mare"
$ws.Range("D67").Value = "Interes Serviciu"
$ws.Range("B69").Value = 30
$ws.Range("C69").Value = "Acasa-Birou"
$ws.Range("D69").Value = " "
$ws.Range("B71").Value = 152
$ws.Range("C71").Value = "Cluj-Cmp. Turzii"
$ws.Range("D71").Value = "Interes Serviciu"
$ws.Range("B73").Value = 356
$ws.Range("C73").Value = "Cluj-Baia-Mare"
$ws.Range("D73").Value = "Interes Serviciu"
$ws.Range("B75").Value = 152
$ws.Range("C75").Value = "Cluj-Cmp. Turzii"
$ws.Range("D75").Value = "Interes Serviciu"

# Totals
$ws.Range("B76").Value = 5415
$ws.Range("B77").Value = 57790

# Signature line with updated submission date
$ws.Range("A87").Value = "Semnătură utilizator:`t`t`t  Data predarii: 17.04.2022"
